# chore: fix batch templates
#
# Replace the placeholder/redacted sample rows (XXXXXX-style NRC/phone
# numbers, flat 20000 "Amount") with realistic sample data, and add a new
# "AccountType" column (H) describing mobile/bank accounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Number formats -------------------------------------------------
# MobileNumber (D) and AccountNumber (F) are identifier-like strings -
# format as Text so leading zeros etc. are preserved. Amount (G) gets a
# 2-decimal numeric format. New AccountType column (H) is Text as well.
$ws.Columns("D").NumberFormat = "@"
$ws.Columns("F").NumberFormat = "@"
$ws.Columns("G").NumberFormat = "0.00"
$ws.Columns("H").NumberFormat = "@"

# --- Header row -------------------------------------------------------
$ws.Range("H1").Value = "AccountType"

# --- Row 2 (Zenaih Tasintha) -------------------------------------------
$ws.Range("D2").Value = "0965123456"
$ws.Range("E2").Value = "334982/10/1"
$ws.Range("F2").Value = "0967123456"
$ws.Range("G2").Value = 10.12
$ws.Range("H2").Value = "mobile"

# --- Row 3 (Hope Tisungeni) --------------------------------------------
$ws.Range("A3").Value = "Hope"
$ws.Range("B3").Value = "Tisungeni"
$ws.Range("C3").Value = "ht@example.com"
$ws.Range("D3").Value = "0975123456"
$ws.Range("E3").Value = "123456/64/1"
$ws.Range("F3").Value = "1234567890123"
$ws.Range("G3").Value = 100
$ws.Range("H3").Value = "bank"

# --- Row 4 (Gift Tumone) ------------------------------------------------
$ws.Range("A4").Value = "Gift"
$ws.Range("B4").Value = "Tumone"
$ws.Range("D4").Value = "0951234567"
$ws.Range("E4").Value = "123456/64/1"
$ws.Range("F4").Value = "0972123321"
$ws.Range("G4").Value = 50.54
$ws.Range("H4").Value = "mobile"

# --- Selection ----------------------------------------------------------
$ws.Range("F3").Select() | Out-Null
